# Auto-generated Excel COM-interop edit script
# Applies the numeric precision corrections described in the commit:
# "bug fix in pf result creation; updated pf result files"
$wb = $excel.ActiveWorkbook

# Sheet: LL_max_6
$ws = $wb.Worksheets.Item("LL_max_6")
$ws.Range("O3").Value = 0.2917056037512947
$ws.Range("P3").Value = 0.8906876332644377
$ws.Range("Q3").Value = 7.091324389062582
$ws.Range("C4").Value = 3.709154423937076
$ws.Range("D4").Value = 3.709154423937076
$ws.Range("F4").Value = 42.82962610251925
$ws.Range("G4").Value = 42.82962610251925
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.9526279648032088
$ws.Range("P4").Value = 0.952627964804755
$ws.Range("Q4").Value = -0.00000000001246565770098904
$ws.Range("N5").Value = 0.9526279647942851
$ws.Range("P5").Value = 0.9526279648136786
$ws.Range("Q5").Value = 0.0000000004605760669781345
$ws.Range("P6").Value = 0.9526279648136786
$ws.Range("Q6").Value = 0.0000000004605760669781345

# Sheet: LL_max_10
$ws = $wb.Worksheets.Item("LL_max_10")
$ws.Range("O3").Value = 0.2917056037512947
$ws.Range("P3").Value = 0.8906876332644377
$ws.Range("Q3").Value = 7.091324389062582
$ws.Range("C4").Value = 3.709154423937076
$ws.Range("D4").Value = 3.709154423937076
$ws.Range("F4").Value = 42.82962610251925
$ws.Range("G4").Value = 42.82962610251925
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.9526279648032088
$ws.Range("P4").Value = 0.952627964804755
$ws.Range("Q4").Value = -0.00000000001246565770098904
$ws.Range("N5").Value = 0.9526279647942851
$ws.Range("P5").Value = 0.9526279648136786
$ws.Range("Q5").Value = 0.0000000004605760669781345
$ws.Range("P6").Value = 0.9526279648136786
$ws.Range("Q6").Value = 0.0000000004605760669781345

# Sheet: LL_max_fault_6
$ws = $wb.Worksheets.Item("LL_max_fault_6")
$ws.Range("O3").Value = 0.2917056037512947
$ws.Range("P3").Value = 0.8906876332644377
$ws.Range("Q3").Value = 7.091324389062582
$ws.Range("C4").Value = 3.709154423937076
$ws.Range("D4").Value = 3.709154423937076
$ws.Range("F4").Value = 42.82962610251925
$ws.Range("G4").Value = 42.82962610251925
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.9526279648032088
$ws.Range("P4").Value = 0.952627964804755
$ws.Range("Q4").Value = -0.00000000001246565770098904
$ws.Range("N5").Value = 0.9526279647942851
$ws.Range("P5").Value = 0.9526279648136786
$ws.Range("Q5").Value = 0.0000000004605760669781345
$ws.Range("P6").Value = 0.9526279648136786
$ws.Range("Q6").Value = 0.0000000004605760669781345

# Sheet: LL_max_fault_10
$ws = $wb.Worksheets.Item("LL_max_fault_10")
$ws.Range("O3").Value = 0.2917056037512947
$ws.Range("P3").Value = 0.8906876332644377
$ws.Range("Q3").Value = 7.091324389062582
$ws.Range("C4").Value = 3.709154423937076
$ws.Range("D4").Value = 3.709154423937076
$ws.Range("F4").Value = 42.82962610251925
$ws.Range("G4").Value = 42.82962610251925
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.9526279648032088
$ws.Range("P4").Value = 0.952627964804755
$ws.Range("Q4").Value = -0.00000000001246565770098904
$ws.Range("N5").Value = 0.9526279647942851
$ws.Range("P5").Value = 0.9526279648136786
$ws.Range("Q5").Value = 0.0000000004605760669781345
$ws.Range("P6").Value = 0.9526279648136786
$ws.Range("Q6").Value = 0.0000000004605760669781345

# Sheet: LL_min_6
$ws = $wb.Worksheets.Item("LL_min_6")
$ws.Range("N2").Value = 0.9821968591989457
$ws.Range("N3").Value = 0.9900703401370795
$ws.Range("P3").Value = 0.7649715466411117
$ws.Range("Q3").Value = 8.150888684274463
$ws.Range("R3").Value = -129.0915432241894
$ws.Range("C4").Value = 3.11908582551553
$ws.Range("D4").Value = 3.11908582551553
$ws.Range("F4").Value = 36.01610081973875
$ws.Range("G4").Value = 36.01610081973875
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.8660254037880062
$ws.Range("P4").Value = 0.8660254037883346
$ws.Range("Q4").Value = 0.0000000001723231751794642
$ws.Range("N5").Value = 0.8660254037798936
$ws.Range("P5").Value = 0.8660254037964472
$ws.Range("Q5").Value = 0.000000001080561154237762
$ws.Range("N6").Value = 0.8660254037798936
$ws.Range("P6").Value = 0.8660254037964472
$ws.Range("Q6").Value = 0.000000001080561154237762

# Sheet: LL_min_10
$ws = $wb.Worksheets.Item("LL_min_10")
$ws.Range("N2").Value = 0.9821968591989457
$ws.Range("N3").Value = 0.9900703401370795
$ws.Range("P3").Value = 0.7649715466411117
$ws.Range("Q3").Value = 8.150888684274463
$ws.Range("R3").Value = -129.0915432241894
$ws.Range("C4").Value = 3.11908582551553
$ws.Range("D4").Value = 3.11908582551553
$ws.Range("F4").Value = 36.01610081973875
$ws.Range("G4").Value = 36.01610081973875
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.8660254037880062
$ws.Range("P4").Value = 0.8660254037883346
$ws.Range("Q4").Value = 0.0000000001723231751794642
$ws.Range("N5").Value = 0.8660254037798936
$ws.Range("P5").Value = 0.8660254037964472
$ws.Range("Q5").Value = 0.000000001080561154237762
$ws.Range("N6").Value = 0.8660254037798936
$ws.Range("P6").Value = 0.8660254037964472
$ws.Range("Q6").Value = 0.000000001080561154237762

# Sheet: LL_min_fault_6
$ws = $wb.Worksheets.Item("LL_min_fault_6")
$ws.Range("N2").Value = 0.9821968591989457
$ws.Range("N3").Value = 0.9900703401370795
$ws.Range("P3").Value = 0.7649715466411117
$ws.Range("Q3").Value = 8.150888684274463
$ws.Range("R3").Value = -129.0915432241894
$ws.Range("C4").Value = 3.11908582551553
$ws.Range("D4").Value = 3.11908582551553
$ws.Range("F4").Value = 36.01610081973875
$ws.Range("G4").Value = 36.01610081973875
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.8660254037880062
$ws.Range("P4").Value = 0.8660254037883346
$ws.Range("Q4").Value = 0.0000000001723231751794642
$ws.Range("N5").Value = 0.8660254037798936
$ws.Range("P5").Value = 0.8660254037964472
$ws.Range("Q5").Value = 0.000000001080561154237762
$ws.Range("N6").Value = 0.8660254037798936
$ws.Range("P6").Value = 0.8660254037964472
$ws.Range("Q6").Value = 0.000000001080561154237762

# Sheet: LL_min_fault_10
$ws = $wb.Worksheets.Item("LL_min_fault_10")
$ws.Range("N2").Value = 0.9821968591989457
$ws.Range("N3").Value = 0.9900703401370795
$ws.Range("P3").Value = 0.7649715466411117
$ws.Range("Q3").Value = 8.150888684274463
$ws.Range("R3").Value = -129.0915432241894
$ws.Range("C4").Value = 3.11908582551553
$ws.Range("D4").Value = 3.11908582551553
$ws.Range("F4").Value = 36.01610081973875
$ws.Range("G4").Value = 36.01610081973875
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.8660254037880062
$ws.Range("P4").Value = 0.8660254037883346
$ws.Range("Q4").Value = 0.0000000001723231751794642
$ws.Range("N5").Value = 0.8660254037798936
$ws.Range("P5").Value = 0.8660254037964472
$ws.Range("Q5").Value = 0.000000001080561154237762
$ws.Range("N6").Value = 0.8660254037798936
$ws.Range("P6").Value = 0.8660254037964472
$ws.Range("Q6").Value = 0.000000001080561154237762

# Sheet: LG_max_6
$ws = $wb.Worksheets.Item("LG_max_6")
$ws.Range("N2").Value = 1.016339946986808
$ws.Range("O2").Value = 1.10000002384591
$ws.Range("P2").Value = 1.037926344349505
$ws.Range("Q2").Value = 31.42118057405401
$ws.Range("R2").Value = -89.99999999999631
$ws.Range("S2").Value = 146.6794423544454
$ws.Range("N3").Value = 0.6615759752774873
$ws.Range("O3").Value = 1.100000023848711
$ws.Range("P3").Value = 0.8500122842528586
$ws.Range("Q3").Value = 39.46791615537077
$ws.Range("R3").Value = -89.99999999999679
$ws.Range("S3").Value = 126.9302900522396
$ws.Range("B4").Value = 3.985202435505792
$ws.Range("E4").Value = 46.01715397828843
$ws.Range("H4").Value = 1.857796819840955
$ws.Range("I4").Value = 3.243741825190523
$ws.Range("J4").Value = 0.6277319162689027
$ws.Range("K4").Value = 2.898438697960071
$ws.Range("L4").Value = 0.6277319163639522
$ws.Range("M4").Value = 2.898438697959879
$ws.Range("N4").Value = 0.6036897740638559
$ws.Range("O4").Value = 1.100000023848711
$ws.Range("P4").Value = 0.7196945219648978
$ws.Range("Q4").Value = 52.69992274962625
$ws.Range("R4").Value = -89.99999999999686
$ws.Range("S4").Value = 120.5514846280052
$ws.Range("T4").Value = 3.985202435505792
$ws.Range("N5").Value = 0.6036897740766251
$ws.Range("O5").Value = 1.100000023848711
$ws.Range("P5").Value = 0.7196945219595307
$ws.Range("Q5").Value = 52.69992274968281
$ws.Range("R5").Value = -89.99999999999699
$ws.Range("S5").Value = 120.5514846289287
$ws.Range("N6").Value = 0.6036897740766251
$ws.Range("O6").Value = 1.100000023848711
$ws.Range("P6").Value = 0.7196945219595307
$ws.Range("Q6").Value = 52.69992274968281
$ws.Range("R6").Value = -89.99999999999699
$ws.Range("S6").Value = 120.5514846289287

# Sheet: LG_max_10
$ws = $wb.Worksheets.Item("LG_max_10")
$ws.Range("N2").Value = 1.016339946986808
$ws.Range("O2").Value = 1.10000002384591
$ws.Range("P2").Value = 1.037926344349505
$ws.Range("Q2").Value = 31.42118057405401
$ws.Range("R2").Value = -89.99999999999631
$ws.Range("S2").Value = 146.6794423544454
$ws.Range("N3").Value = 0.6615759752774873
$ws.Range("O3").Value = 1.100000023848711
$ws.Range("P3").Value = 0.8500122842528586
$ws.Range("Q3").Value = 39.46791615537077
$ws.Range("R3").Value = -89.99999999999679
$ws.Range("S3").Value = 126.9302900522396
$ws.Range("B4").Value = 3.985202435505792
$ws.Range("E4").Value = 46.01715397828843
$ws.Range("H4").Value = 1.857796819840955
$ws.Range("I4").Value = 3.243741825190523
$ws.Range("J4").Value = 0.6277319162689027
$ws.Range("K4").Value = 2.898438697960071
$ws.Range("L4").Value = 0.6277319163639522
$ws.Range("M4").Value = 2.898438697959879
$ws.Range("N4").Value = 0.6036897740638559
$ws.Range("O4").Value = 1.100000023848711
$ws.Range("P4").Value = 0.7196945219648978
$ws.Range("Q4").Value = 52.69992274962625
$ws.Range("R4").Value = -89.99999999999686
$ws.Range("S4").Value = 120.5514846280052
$ws.Range("T4").Value = 3.985202435505792
$ws.Range("N5").Value = 0.6036897740766251
$ws.Range("O5").Value = 1.100000023848711
$ws.Range("P5").Value = 0.7196945219595307
$ws.Range("Q5").Value = 52.69992274968281
$ws.Range("R5").Value = -89.99999999999699
$ws.Range("S5").Value = 120.5514846289287
$ws.Range("N6").Value = 0.6036897740766251
$ws.Range("O6").Value = 1.100000023848711
$ws.Range("P6").Value = 0.7196945219595307
$ws.Range("Q6").Value = 52.69992274968281
$ws.Range("R6").Value = -89.99999999999699
$ws.Range("S6").Value = 120.5514846289287

# Sheet: LG_max_fault_6
$ws = $wb.Worksheets.Item("LG_max_fault_6")
$ws.Range("N2").Value = 1.072156783213623
$ws.Range("O2").Value = 1.100000023843317
$ws.Range("P2").Value = 1.086893850314861
$ws.Range("Q2").Value = 29.96661352520651
$ws.Range("R2").Value = -89.99999999999635
$ws.Range("S2").Value = 148.7122982625788
$ws.Range("N3").Value = 0.9542051367681318
$ws.Range("O3").Value = 1.100000023844832
$ws.Range("P3").Value = 1.044970925665078
$ws.Range("Q3").Value = 29.33776687354478
$ws.Range("R3").Value = -89.99999999999659
$ws.Range("S3").Value = 142.7523590701965
$ws.Range("B4").Value = 1.265928718513589
$ws.Range("E4").Value = 14.61768572817398
$ws.Range("H4").Value = 1.857796819840955
$ws.Range("I4").Value = 3.243741825190523
$ws.Range("J4").Value = 0.6277319162689027
$ws.Range("K4").Value = 2.898438697960071
$ws.Range("L4").Value = 0.6277319163639522
$ws.Range("M4").Value = 2.898438697959879
$ws.Range("N4").Value = 0.9132255744650989
$ws.Range("O4").Value = 1.100000023844832
$ws.Range("P4").Value = 1.003516002765787
$ws.Range("Q4").Value = 31.07240830832294
$ws.Range("R4").Value = -89.9999999999966
$ws.Range("S4").Value = 141.2103841214374
$ws.Range("T4").Value = 1.265928718513589
$ws.Range("N5").Value = 0.9132255744681027
$ws.Range("O5").Value = 1.100000023844832
$ws.Range("P5").Value = 1.003516002764256
$ws.Range("Q5").Value = 31.07240830849433
$ws.Range("R5").Value = -89.99999999999666
$ws.Range("S5").Value = 141.2103841216521
$ws.Range("N6").Value = 0.9132255744681027
$ws.Range("O6").Value = 1.100000023844832
$ws.Range("P6").Value = 1.003516002764256
$ws.Range("Q6").Value = 31.07240830849433
$ws.Range("R6").Value = -89.99999999999666
$ws.Range("S6").Value = 141.2103841216521

# Sheet: LG_max_fault_10
$ws = $wb.Worksheets.Item("LG_max_fault_10")
$ws.Range("N2").Value = 1.072156783213623
$ws.Range("O2").Value = 1.100000023843317
$ws.Range("P2").Value = 1.086893850314861
$ws.Range("Q2").Value = 29.96661352520651
$ws.Range("R2").Value = -89.99999999999635
$ws.Range("S2").Value = 148.7122982625788
$ws.Range("N3").Value = 0.9542051367681318
$ws.Range("O3").Value = 1.100000023844832
$ws.Range("P3").Value = 1.044970925665078
$ws.Range("Q3").Value = 29.33776687354478
$ws.Range("R3").Value = -89.99999999999659
$ws.Range("S3").Value = 142.7523590701965
$ws.Range("B4").Value = 1.265928718513589
$ws.Range("E4").Value = 14.61768572817398
$ws.Range("H4").Value = 1.857796819840955
$ws.Range("I4").Value = 3.243741825190523
$ws.Range("J4").Value = 0.6277319162689027
$ws.Range("K4").Value = 2.898438697960071
$ws.Range("L4").Value = 0.6277319163639522
$ws.Range("M4").Value = 2.898438697959879
$ws.Range("N4").Value = 0.9132255744650989
$ws.Range("O4").Value = 1.100000023844832
$ws.Range("P4").Value = 1.003516002765787
$ws.Range("Q4").Value = 31.07240830832294
$ws.Range("R4").Value = -89.9999999999966
$ws.Range("S4").Value = 141.2103841214374
$ws.Range("T4").Value = 1.265928718513589
$ws.Range("N5").Value = 0.9132255744681027
$ws.Range("O5").Value = 1.100000023844832
$ws.Range("P5").Value = 1.003516002764256
$ws.Range("Q5").Value = 31.07240830849433
$ws.Range("R5").Value = -89.99999999999666
$ws.Range("S5").Value = 141.2103841216521
$ws.Range("N6").Value = 0.9132255744681027
$ws.Range("O6").Value = 1.100000023844832
$ws.Range("P6").Value = 1.003516002764256
$ws.Range("Q6").Value = 31.07240830849433
$ws.Range("R6").Value = -89.99999999999666
$ws.Range("S6").Value = 141.2103841216521

# Sheet: LG_min_6
$ws = $wb.Worksheets.Item("LG_min_6")
$ws.Range("N2").Value = 0.9208231165114505
$ws.Range("O2").Value = 0.9999999999969289
$ws.Range("P2").Value = 0.9583183403729527
$ws.Range("Q2").Value = 30.31368549467228
$ws.Range("R2").Value = -89.99999999999633
$ws.Range("S2").Value = 146.0472524790669
$ws.Range("N3").Value = 0.6161181908425291
$ws.Range("O3").Value = 1.000000000001022
$ws.Range("P3").Value = 0.8693567795104428
$ws.Range("Q3").Value = 30.41440561874932
$ws.Range("R3").Value = -89.999999999997
$ws.Range("S3").Value = 127.6747470341949
$ws.Range("B4").Value = 3.174039680767477
$ws.Range("E4").Value = 36.65065328219313
$ws.Range("H4").Value = 3.566942654999328
$ws.Range("I4").Value = 3.235808697828856
$ws.Range("J4").Value = 1.113751880275152
$ws.Range("K4").Value = 3.006397543745624
$ws.Range("L4").Value = 1.113751880233234
$ws.Range("N4").Value = 0.5408411981898492
$ws.Range("O4").Value = 1.000000000001022
$ws.Range("P4").Value = 0.7068031750084121
$ws.Range("Q4").Value = 47.14357112007056
$ws.Range("R4").Value = -89.99999999999707
$ws.Range("S4").Value = 121.3628957621809
$ws.Range("T4").Value = 3.174039680767477
$ws.Range("N5").Value = 0.5408411982060234
$ws.Range("O5").Value = 1.000000000001022
$ws.Range("P5").Value = 0.7068031750053295
$ws.Range("Q5").Value = 47.14357111992546
$ws.Range("R5").Value = -89.99999999999729
$ws.Range("S5").Value = 121.3628957634726
$ws.Range("N6").Value = 0.5408411982060234
$ws.Range("O6").Value = 1.000000000001022
$ws.Range("P6").Value = 0.7068031750053295
$ws.Range("Q6").Value = 47.14357111992546
$ws.Range("R6").Value = -89.99999999999729
$ws.Range("S6").Value = 121.3628957634726

# Sheet: LG_min_10
$ws = $wb.Worksheets.Item("LG_min_10")
$ws.Range("N2").Value = 0.9208231165114505
$ws.Range("O2").Value = 0.9999999999969289
$ws.Range("P2").Value = 0.9583183403729527
$ws.Range("Q2").Value = 30.31368549467228
$ws.Range("R2").Value = -89.99999999999633
$ws.Range("S2").Value = 146.0472524790669
$ws.Range("N3").Value = 0.6161181908425291
$ws.Range("O3").Value = 1.000000000001022
$ws.Range("P3").Value = 0.8693567795104428
$ws.Range("Q3").Value = 30.41440561874932
$ws.Range("R3").Value = -89.999999999997
$ws.Range("S3").Value = 127.6747470341949
$ws.Range("B4").Value = 3.174039680767477
$ws.Range("E4").Value = 36.65065328219313
$ws.Range("H4").Value = 3.566942654999328
$ws.Range("I4").Value = 3.235808697828856
$ws.Range("J4").Value = 1.113751880275152
$ws.Range("K4").Value = 3.006397543745624
$ws.Range("L4").Value = 1.113751880233234
$ws.Range("N4").Value = 0.5408411981898492
$ws.Range("O4").Value = 1.000000000001022
$ws.Range("P4").Value = 0.7068031750084121
$ws.Range("Q4").Value = 47.14357112007056
$ws.Range("R4").Value = -89.99999999999707
$ws.Range("S4").Value = 121.3628957621809
$ws.Range("T4").Value = 3.174039680767477
$ws.Range("N5").Value = 0.5408411982060234
$ws.Range("O5").Value = 1.000000000001022
$ws.Range("P5").Value = 0.7068031750053295
$ws.Range("Q5").Value = 47.14357111992546
$ws.Range("R5").Value = -89.99999999999729
$ws.Range("S5").Value = 121.3628957634726
$ws.Range("N6").Value = 0.5408411982060234
$ws.Range("O6").Value = 1.000000000001022
$ws.Range("P6").Value = 0.7068031750053295
$ws.Range("Q6").Value = 47.14357111992546
$ws.Range("R6").Value = -89.99999999999729
$ws.Range("S6").Value = 121.3628957634726

# Sheet: LG_min_fault_6
$ws = $wb.Worksheets.Item("LG_min_fault_6")
$ws.Range("N2").Value = 0.9729910730890348
$ws.Range("O2").Value = 0.9999999999987059
$ws.Range("P2").Value = 0.9888062153124251
$ws.Range("Q2").Value = 29.86351783395289
$ws.Range("S2").Value = 148.5772221076714
$ws.Range("N3").Value = 0.870254229744171
$ws.Range("O3").Value = 1.000000000000083
$ws.Range("P3").Value = 0.958347064291428
$ws.Range("Q3").Value = 28.81566541771512
$ws.Range("R3").Value = -89.99999999999662
$ws.Range("S3").Value = 142.7153207693948
$ws.Range("B4").Value = 1.084437662235425
$ws.Range("E4").Value = 12.52200752421982
$ws.Range("H4").Value = 3.566942654999328
$ws.Range("I4").Value = 3.235808697828856
$ws.Range("J4").Value = 1.113751880275152
$ws.Range("K4").Value = 3.006397543745624
$ws.Range("L4").Value = 1.113751880233234
$ws.Range("N4").Value = 0.8291717042090867
$ws.Range("O4").Value = 1.000000000000083
$ws.Range("P4").Value = 0.8989295800157858
$ws.Range("Q4").Value = 32.02702342589205
$ws.Range("R4").Value = -89.99999999999666
$ws.Range("S4").Value = 141.444800011689
$ws.Range("T4").Value = 1.084437662235425
$ws.Range("N5").Value = 0.8291717042143361
$ws.Range("O5").Value = 1.000000000000083
$ws.Range("P5").Value = 0.8989295800158428
$ws.Range("Q5").Value = 32.02702342601567
$ws.Range("R5").Value = -89.99999999999673
$ws.Range("S5").Value = 141.4448000120424
$ws.Range("N6").Value = 0.8291717042143361
$ws.Range("O6").Value = 1.000000000000083
$ws.Range("P6").Value = 0.8989295800158428
$ws.Range("Q6").Value = 32.02702342601567
$ws.Range("R6").Value = -89.99999999999673
$ws.Range("S6").Value = 141.4448000120424

# Sheet: LG_min_fault_10
$ws = $wb.Worksheets.Item("LG_min_fault_10")
$ws.Range("N2").Value = 0.9729910730890348
$ws.Range("O2").Value = 0.9999999999987059
$ws.Range("P2").Value = 0.9888062153124251
$ws.Range("Q2").Value = 29.86351783395289
$ws.Range("S2").Value = 148.5772221076714
$ws.Range("N3").Value = 0.870254229744171
$ws.Range("O3").Value = 1.000000000000083
$ws.Range("P3").Value = 0.958347064291428
$ws.Range("Q3").Value = 28.81566541771512
$ws.Range("R3").Value = -89.99999999999662
$ws.Range("S3").Value = 142.7153207693948
$ws.Range("B4").Value = 1.084437662235425
$ws.Range("E4").Value = 12.52200752421982
$ws.Range("H4").Value = 3.566942654999328
$ws.Range("I4").Value = 3.235808697828856
$ws.Range("J4").Value = 1.113751880275152
$ws.Range("K4").Value = 3.006397543745624
$ws.Range("L4").Value = 1.113751880233234
$ws.Range("N4").Value = 0.8291717042090867
$ws.Range("O4").Value = 1.000000000000083
$ws.Range("P4").Value = 0.8989295800157858
$ws.Range("Q4").Value = 32.02702342589205
$ws.Range("R4").Value = -89.99999999999666
$ws.Range("S4").Value = 141.444800011689
$ws.Range("T4").Value = 1.084437662235425
$ws.Range("N5").Value = 0.8291717042143361
$ws.Range("O5").Value = 1.000000000000083
$ws.Range("P5").Value = 0.8989295800158428
$ws.Range("Q5").Value = 32.02702342601567
$ws.Range("R5").Value = -89.99999999999673
$ws.Range("S5").Value = 141.4448000120424
$ws.Range("N6").Value = 0.8291717042143361
$ws.Range("O6").Value = 1.000000000000083
$ws.Range("P6").Value = 0.8989295800158428
$ws.Range("Q6").Value = 32.02702342601567
$ws.Range("R6").Value = -89.99999999999673
$ws.Range("S6").Value = 141.4448000120424

# Sheet: LLG_max_6
$ws = $wb.Worksheets.Item("LLG_max_6")
$ws.Range("N2").Value = 1.02986751035138
$ws.Range("O2").Value = 0.9380319009590135
$ws.Range("P2").Value = 1.025521871529738
$ws.Range("Q2").Value = 26.25928376955642
$ws.Range("R2").Value = -91.13019952197448
$ws.Range("S2").Value = 151.9526811209942
$ws.Range("N3").Value = 0.8279179750199849
$ws.Range("O3").Value = 0.2917056037508406
$ws.Range("P3").Value = 0.7152018015263178
$ws.Range("Q3").Value = 3.258617650687572
$ws.Range("R3").Value = -119.1713773543864
$ws.Range("S3").Value = 163.1221622147212
$ws.Range("C4").Value = 4.487274272340209
$ws.Range("D4").Value = 3.763988065653138
$ws.Range("F4").Value = 51.81458018126604
$ws.Range("G4").Value = 43.46279045862757
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.6886390086961146
$ws.Range("P4").Value = 0.6886390087023841
$ws.Range("Q4").Value = -4.968885503490361
$ws.Range("S4").Value = 175.0311144964618
$ws.Range("T4").Value = 3.684449340743881
$ws.Range("N5").Value = 0.6886390086889234
$ws.Range("P5").Value = 0.6886390087143381
$ws.Range("Q5").Value = -4.96888550245307
$ws.Range("S5").Value = 175.0311144963239
$ws.Range("N6").Value = 0.6886390086889234
$ws.Range("P6").Value = 0.6886390087143381
$ws.Range("Q6").Value = -4.96888550245307
$ws.Range("S6").Value = 175.0311144963239

# Sheet: LLG_max_10
$ws = $wb.Worksheets.Item("LLG_max_10")
$ws.Range("N2").Value = 1.02986751035138
$ws.Range("O2").Value = 0.9380319009590135
$ws.Range("P2").Value = 1.025521871529738
$ws.Range("Q2").Value = 26.25928376955642
$ws.Range("R2").Value = -91.13019952197448
$ws.Range("S2").Value = 151.9526811209942
$ws.Range("N3").Value = 0.8279179750199849
$ws.Range("O3").Value = 0.2917056037508406
$ws.Range("P3").Value = 0.7152018015263178
$ws.Range("Q3").Value = 3.258617650687572
$ws.Range("R3").Value = -119.1713773543864
$ws.Range("S3").Value = 163.1221622147212
$ws.Range("C4").Value = 4.487274272340209
$ws.Range("D4").Value = 3.763988065653138
$ws.Range("F4").Value = 51.81458018126604
$ws.Range("G4").Value = 43.46279045862757
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.6886390086961146
$ws.Range("P4").Value = 0.6886390087023841
$ws.Range("Q4").Value = -4.968885503490361
$ws.Range("S4").Value = 175.0311144964618
$ws.Range("T4").Value = 3.684449340743881
$ws.Range("N5").Value = 0.6886390086889234
$ws.Range("P5").Value = 0.6886390087143381
$ws.Range("Q5").Value = -4.96888550245307
$ws.Range("S5").Value = 175.0311144963239
$ws.Range("N6").Value = 0.6886390086889234
$ws.Range("P6").Value = 0.6886390087143381
$ws.Range("Q6").Value = -4.96888550245307
$ws.Range("S6").Value = 175.0311144963239

# Sheet: LLG_max_fault_6
$ws = $wb.Worksheets.Item("LLG_max_fault_6")
$ws.Range("N2").Value = 1.02986751035138
$ws.Range("O2").Value = 0.9380319009590135
$ws.Range("P2").Value = 1.025521871529738
$ws.Range("Q2").Value = 26.25928376955642
$ws.Range("R2").Value = -91.13019952197448
$ws.Range("S2").Value = 151.9526811209942
$ws.Range("N3").Value = 0.8279179750199849
$ws.Range("O3").Value = 0.2917056037508406
$ws.Range("P3").Value = 0.7152018015263178
$ws.Range("Q3").Value = 3.258617650687572
$ws.Range("R3").Value = -119.1713773543864
$ws.Range("S3").Value = 163.1221622147212
$ws.Range("C4").Value = 4.487274272340209
$ws.Range("D4").Value = 3.763988065653138
$ws.Range("F4").Value = 51.81458018126604
$ws.Range("G4").Value = 43.46279045862757
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.6886390086961146
$ws.Range("P4").Value = 0.6886390087023841
$ws.Range("Q4").Value = -4.968885503490361
$ws.Range("S4").Value = 175.0311144964618
$ws.Range("T4").Value = 3.684449340743881
$ws.Range("N5").Value = 0.6886390086889234
$ws.Range("P5").Value = 0.6886390087143381
$ws.Range("Q5").Value = -4.96888550245307
$ws.Range("S5").Value = 175.0311144963239
$ws.Range("N6").Value = 0.6886390086889234
$ws.Range("P6").Value = 0.6886390087143381
$ws.Range("Q6").Value = -4.96888550245307
$ws.Range("S6").Value = 175.0311144963239

# Sheet: LLG_max_fault_10
$ws = $wb.Worksheets.Item("LLG_max_fault_10")
$ws.Range("N2").Value = 1.02986751035138
$ws.Range("O2").Value = 0.9380319009590135
$ws.Range("P2").Value = 1.025521871529738
$ws.Range("Q2").Value = 26.25928376955642
$ws.Range("R2").Value = -91.13019952197448
$ws.Range("S2").Value = 151.9526811209942
$ws.Range("N3").Value = 0.8279179750199849
$ws.Range("O3").Value = 0.2917056037508406
$ws.Range("P3").Value = 0.7152018015263178
$ws.Range("Q3").Value = 3.258617650687572
$ws.Range("R3").Value = -119.1713773543864
$ws.Range("S3").Value = 163.1221622147212
$ws.Range("C4").Value = 4.487274272340209
$ws.Range("D4").Value = 3.763988065653138
$ws.Range("F4").Value = 51.81458018126604
$ws.Range("G4").Value = 43.46279045862757
$ws.Range("H4").Value = 1.857796819840932
$ws.Range("I4").Value = 3.24374182519076
$ws.Range("J4").Value = 0.6277319163569953
$ws.Range("K4").Value = 2.89843869794806
$ws.Range("L4").Value = 0.6277319163639853
$ws.Range("M4").Value = 2.898438697959971
$ws.Range("N4").Value = 0.6886390086961146
$ws.Range("P4").Value = 0.6886390087023841
$ws.Range("Q4").Value = -4.968885503490361
$ws.Range("S4").Value = 175.0311144964618
$ws.Range("T4").Value = 3.684449340743881
$ws.Range("N5").Value = 0.6886390086889234
$ws.Range("P5").Value = 0.6886390087143381
$ws.Range("Q5").Value = -4.96888550245307
$ws.Range("S5").Value = 175.0311144963239
$ws.Range("N6").Value = 0.6886390086889234
$ws.Range("P6").Value = 0.6886390087143381
$ws.Range("Q6").Value = -4.96888550245307
$ws.Range("S6").Value = 175.0311144963239

# Sheet: LLG_min_6
$ws = $wb.Worksheets.Item("LLG_min_6")
$ws.Range("N2").Value = 0.9482439689260066
$ws.Range("O2").Value = 0.8500049525399591
$ws.Range("P2").Value = 0.9311545463007321
$ws.Range("Q2").Value = 25.25873620967787
$ws.Range("R2").Value = -92.65381631332932
$ws.Range("S2").Value = 151.4886012729938
$ws.Range("N3").Value = 0.8500087908957914
$ws.Range("O3").Value = 0.3617201835042794
$ws.Range("P3").Value = 0.6672173733277518
$ws.Range("Q3").Value = 2.497096978932195
$ws.Range("R3").Value = -129.0915432241547
$ws.Range("S3").Value = 158.5760776527842
$ws.Range("C4").Value = 3.830038527961571
$ws.Range("D4").Value = 2.93072060203517
$ws.Range("F4").Value = 44.22547550250501
$ws.Range("G4").Value = 33.84104657009175
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.6639160778413813
$ws.Range("P4").Value = 0.6639160778412487
$ws.Range("Q4").Value = -6.75741831289674
$ws.Range("S4").Value = 173.2425816865416
$ws.Range("T4").Value = 2.757145480035794
$ws.Range("N5").Value = 0.6639160778354332
$ws.Range("P5").Value = 0.6639160778546471
$ws.Range("Q5").Value = -6.757418311286667
$ws.Range("S5").Value = 173.2425816859635
$ws.Range("N6").Value = 0.6639160778354332
$ws.Range("P6").Value = 0.6639160778546471
$ws.Range("Q6").Value = -6.757418311286667
$ws.Range("S6").Value = 173.2425816859635

# Sheet: LLG_min_10
$ws = $wb.Worksheets.Item("LLG_min_10")
$ws.Range("N2").Value = 0.9482439689260066
$ws.Range("O2").Value = 0.8500049525399591
$ws.Range("P2").Value = 0.9311545463007321
$ws.Range("Q2").Value = 25.25873620967787
$ws.Range("R2").Value = -92.65381631332932
$ws.Range("S2").Value = 151.4886012729938
$ws.Range("N3").Value = 0.8500087908957914
$ws.Range("O3").Value = 0.3617201835042794
$ws.Range("P3").Value = 0.6672173733277518
$ws.Range("Q3").Value = 2.497096978932195
$ws.Range("R3").Value = -129.0915432241547
$ws.Range("S3").Value = 158.5760776527842
$ws.Range("C4").Value = 3.830038527961571
$ws.Range("D4").Value = 2.93072060203517
$ws.Range("F4").Value = 44.22547550250501
$ws.Range("G4").Value = 33.84104657009175
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.6639160778413813
$ws.Range("P4").Value = 0.6639160778412487
$ws.Range("Q4").Value = -6.75741831289674
$ws.Range("S4").Value = 173.2425816865416
$ws.Range("T4").Value = 2.757145480035794
$ws.Range("N5").Value = 0.6639160778354332
$ws.Range("P5").Value = 0.6639160778546471
$ws.Range("Q5").Value = -6.757418311286667
$ws.Range("S5").Value = 173.2425816859635
$ws.Range("N6").Value = 0.6639160778354332
$ws.Range("P6").Value = 0.6639160778546471
$ws.Range("Q6").Value = -6.757418311286667
$ws.Range("S6").Value = 173.2425816859635

# Sheet: LLG_min_fault_6
$ws = $wb.Worksheets.Item("LLG_min_fault_6")
$ws.Range("N2").Value = 0.9482439689260066
$ws.Range("O2").Value = 0.8500049525399591
$ws.Range("P2").Value = 0.9311545463007321
$ws.Range("Q2").Value = 25.25873620967787
$ws.Range("R2").Value = -92.65381631332932
$ws.Range("S2").Value = 151.4886012729938
$ws.Range("N3").Value = 0.8500087908957914
$ws.Range("O3").Value = 0.3617201835042794
$ws.Range("P3").Value = 0.6672173733277518
$ws.Range("Q3").Value = 2.497096978932195
$ws.Range("R3").Value = -129.0915432241547
$ws.Range("S3").Value = 158.5760776527842
$ws.Range("C4").Value = 3.830038527961571
$ws.Range("D4").Value = 2.93072060203517
$ws.Range("F4").Value = 44.22547550250501
$ws.Range("G4").Value = 33.84104657009175
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.6639160778413813
$ws.Range("P4").Value = 0.6639160778412487
$ws.Range("Q4").Value = -6.75741831289674
$ws.Range("S4").Value = 173.2425816865416
$ws.Range("T4").Value = 2.757145480035794
$ws.Range("N5").Value = 0.6639160778354332
$ws.Range("P5").Value = 0.6639160778546471
$ws.Range("Q5").Value = -6.757418311286667
$ws.Range("S5").Value = 173.2425816859635
$ws.Range("N6").Value = 0.6639160778354332
$ws.Range("P6").Value = 0.6639160778546471
$ws.Range("Q6").Value = -6.757418311286667
$ws.Range("S6").Value = 173.2425816859635

# Sheet: LLG_min_fault_10
$ws = $wb.Worksheets.Item("LLG_min_fault_10")
$ws.Range("N2").Value = 0.9482439689260066
$ws.Range("O2").Value = 0.8500049525399591
$ws.Range("P2").Value = 0.9311545463007321
$ws.Range("Q2").Value = 25.25873620967787
$ws.Range("R2").Value = -92.65381631332932
$ws.Range("S2").Value = 151.4886012729938
$ws.Range("N3").Value = 0.8500087908957914
$ws.Range("O3").Value = 0.3617201835042794
$ws.Range("P3").Value = 0.6672173733277518
$ws.Range("Q3").Value = 2.497096978932195
$ws.Range("R3").Value = -129.0915432241547
$ws.Range("S3").Value = 158.5760776527842
$ws.Range("C4").Value = 3.830038527961571
$ws.Range("D4").Value = 2.93072060203517
$ws.Range("F4").Value = 44.22547550250501
$ws.Range("G4").Value = 33.84104657009175
$ws.Range("H4").Value = 3.566942654999266
$ws.Range("I4").Value = 3.235808697829114
$ws.Range("J4").Value = 1.113751880246062
$ws.Range("K4").Value = 3.00639754376794
$ws.Range("L4").Value = 1.113751880233277
$ws.Range("M4").Value = 3.006397543745544
$ws.Range("N4").Value = 0.6639160778413813
$ws.Range("P4").Value = 0.6639160778412487
$ws.Range("Q4").Value = -6.75741831289674
$ws.Range("S4").Value = 173.2425816865416
$ws.Range("T4").Value = 2.757145480035794
$ws.Range("N5").Value = 0.6639160778354332
$ws.Range("P5").Value = 0.6639160778546471
$ws.Range("Q5").Value = -6.757418311286667
$ws.Range("S5").Value = 173.2425816859635
$ws.Range("N6").Value = 0.6639160778354332
$ws.Range("P6").Value = 0.6639160778546471
$ws.Range("Q6").Value = -6.757418311286667
$ws.Range("S6").Value = 173.2425816859635
